$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cells: "Nr." -> "Nr" and "Facultatea" -> "Facultate"
$ws.Range("A1").Value = "Nr"
$ws.Range("C1").Value = "Facultate"

# Move the active selection to C1 (matches the saved selection in the diff)
$ws.Range("C1").Select()
